$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the last existing data row (660) down through the new rows
$ws.Range("A660:G660").Copy()
$ws.Range("A661:G672").PasteSpecial(-4122)

# Populate the new rows with Date, Open, High, Low, Close, AdjClose, Volume
$ws.Range("A661").Value = 45129
$ws.Range("B661").Value = 4550.16015625
$ws.Range("C661").Value = 4555
$ws.Range("D661").Value = 4535.7900390625
$ws.Range("E661").Value = 4536.33984375
$ws.Range("F661").Value = 4536.33984375
$ws.Range("G661").Value = 3570190000

$ws.Range("A662").Value = 45130
$ws.Range("B662").Value = 4550.16015625
$ws.Range("C662").Value = 4555
$ws.Range("D662").Value = 4535.7900390625
$ws.Range("E662").Value = 4536.33984375
$ws.Range("F662").Value = 4536.33984375
$ws.Range("G662").Value = 3570190000

$ws.Range("A663").Value = 45131
$ws.Range("B663").Value = 4543.39013671875
$ws.Range("C663").Value = 4563.41015625
$ws.Range("D663").Value = 4541.2900390625
$ws.Range("E663").Value = 4554.64013671875
$ws.Range("F663").Value = 4554.64013671875
$ws.Range("G663").Value = 3856250000

$ws.Range("A664").Value = 45132
$ws.Range("B664").Value = 4555.18994140625
$ws.Range("C664").Value = 4580.6201171875
$ws.Range("D664").Value = 4552.419921875
$ws.Range("E664").Value = 4567.4599609375
$ws.Range("F664").Value = 4567.4599609375
$ws.Range("G664").Value = 3812470000

$ws.Range("A665").Value = 45133
$ws.Range("B665").Value = 4558.9599609375
$ws.Range("C665").Value = 4582.47021484375
$ws.Range("D665").Value = 4547.580078125
$ws.Range("E665").Value = 4566.75
$ws.Range("F665").Value = 4566.75
$ws.Range("G665").Value = 3990290000

$ws.Range("A666").Value = 45134
$ws.Range("B666").Value = 4598.259765625
$ws.Range("C666").Value = 4607.06982421875
$ws.Range("D666").Value = 4528.56005859375
$ws.Range("E666").Value = 4537.41015625
$ws.Range("F666").Value = 4537.41015625
$ws.Range("G666").Value = 4553210000

$ws.Range("A667").Value = 45135
$ws.Range("B667").Value = 4565.75
$ws.Range("C667").Value = 4590.16015625
$ws.Range("D667").Value = 4564.009765625
$ws.Range("E667").Value = 4582.22998046875
$ws.Range("F667").Value = 4582.22998046875
$ws.Range("G667").Value = 3981010000

$ws.Range("A668").Value = 45136
$ws.Range("B668").Value = 4565.75
$ws.Range("C668").Value = 4590.16015625
$ws.Range("D668").Value = 4564.009765625
$ws.Range("E668").Value = 4582.22998046875
$ws.Range("F668").Value = 4582.22998046875
$ws.Range("G668").Value = 3981010000

$ws.Range("A669").Value = 45137
$ws.Range("B669").Value = 4565.75
$ws.Range("C669").Value = 4590.16015625
$ws.Range("D669").Value = 4564.009765625
$ws.Range("E669").Value = 4582.22998046875
$ws.Range("F669").Value = 4582.22998046875
$ws.Range("G669").Value = 3981010000

$ws.Range("A670").Value = 45138
$ws.Range("B670").Value = 4584.81982421875
$ws.Range("C670").Value = 4594.22021484375
$ws.Range("D670").Value = 4573.14013671875
$ws.Range("E670").Value = 4588.9599609375
$ws.Range("F670").Value = 4588.9599609375
$ws.Range("G670").Value = 4503600000

$ws.Range("A671").Value = 45139
$ws.Range("B671").Value = 4578.830078125
$ws.Range("C671").Value = 4584.6201171875
$ws.Range("D671").Value = 4567.52978515625
$ws.Range("E671").Value = 4576.72998046875
$ws.Range("F671").Value = 4576.72998046875
$ws.Range("G671").Value = 4042370000

$ws.Range("A672").Value = 45140
$ws.Range("B672").Value = 4550.93017578125
$ws.Range("C672").Value = 4550.93017578125
$ws.Range("D672").Value = 4505.75
$ws.Range("E672").Value = 4513.39013671875
$ws.Range("F672").Value = 4513.39013671875
$ws.Range("G672").Value = 4270710000
